$p = $ppt.ActivePresentation

# --- Slide 5: rewrite the "For each ..." bullet (single run, keep its rPr) ---
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(3)
$c5 = $para5.Characters(1, 1)
$c5.Text = "For each components– check special cases –Find upper  &  lower bound(Greedy Algorithm) of the components-run Brute Force"
$tail5 = $para5.Characters(121, 105)
$tail5.Text = ""

# --- Slide 6: "Lower Bound" bullet loses its redundant endParaRPr, and the
#     "Genetic "+"Algorithm" runs merge into a single run ---
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange

# Drop the stray trailing endParaRPr on the "Lower Bound" paragraph by
# deleting it and re-inserting identical text ahead of "Special Cases"
# (which already has the exact rPr/pPr shape we want, with no endParaRPr).
$paraLower = $tr6.Paragraphs(3)
$paraLower.Delete()
$paraSpecial = $tr6.Paragraphs(3)
$paraSpecial.InsertBefore("Lower Bound`r")

# Merge the two "Genetic " / "Algorithm" runs into one run.
$paraGenetic = $tr6.Paragraphs(6)
$cG = $paraGenetic.Characters(1, 8)
$cG.Text = "Genetic Algorithm"
$tailG = $paraGenetic.Characters(18, 9)
$tailG.Text = ""

# --- Slide 9: merge the "LOWER " + "BOUND(GREEDY ALGORITHM)" runs ---
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)
$tr9 = $sh9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(1)
$c9 = $para9.Characters(1, 6)
$c9.Text = "LOWER BOUND(GREEDY ALGORITHM)"
$tail9 = $para9.Characters(30, 23)
$tail9.Text = ""
